# Update the "想去人数" (want-to-go count) values in column F for rows 2-7
# on both the "展览" and "全部类型" worksheets, matching the newly scraped data.

$wb = $excel.ActiveWorkbook

$newValues = @{
    2 = 1306
    3 = 1703
    4 = 67
    5 = 6246
    6 = 90
    7 = 106
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $newValues.Keys) {
        $ws.Range("F$row").Value = $newValues[$row]
    }
}
